$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")
$row = 3

# A3 .. E3
$ws.Cells.Item($row, 1).Value  = 106645464              # A  Id
$ws.Cells.Item($row, 2).Value  = 56369                   # B  Taxonsorteringsordning
$ws.Cells.Item($row, 3).Value  = 'Ovaliderad'             # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value  = 'VU'                     # D  Rodlistade
$ws.Cells.Item($row, 5).Value  = 100136                  # E  TaxonId

# F3 .. I3
$ws.Cells.Item($row, 6).Value  = 'Lappuggla'              # F  Artnamn
$ws.Cells.Item($row, 7).Value  = 'Strix nebulosa'         # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = 'J.R. Forster, 1772'     # H  Auktor

# I3 must stay text ("1"), not be auto-coerced to a number: lead with an
# apostrophe so Excel stores it as a quoted text literal.
$ws.Cells.Item($row, 9).Value  = "'" + '1'                # I  Antal

# K3, L3, N3 are present-but-empty text cells in the source row. A plain
# empty-string assignment is treated as "no content" and the cell is
# dropped, so force the cell to Text format first - that keeps the (blank)
# cell alive in the sheet.
$ws.Cells.Item($row, 11).NumberFormat = "@"
$ws.Cells.Item($row, 11).Value = ""                       # K  Alder-Stadium
$ws.Cells.Item($row, 12).NumberFormat = "@"
$ws.Cells.Item($row, 12).Value = ""                       # L  Kon

$ws.Cells.Item($row, 13).Value = 'lockläte, övriga läten' # M  Aktivitet

$ws.Cells.Item($row, 14).NumberFormat = "@"
$ws.Cells.Item($row, 14).Value = ""                       # N  Metod

$ws.Cells.Item($row, 16).Value = 'Yasjön, Åkulla, Hl'      # P  Lokalnamn
$ws.Cells.Item($row, 17).Value = 351937                   # Q  Ost
$ws.Cells.Item($row, 18).Value = 6334474                  # R  Nord
$ws.Cells.Item($row, 19).Value = 5                        # S  Noggrannhet
$ws.Cells.Item($row, 20).Value = 'Halland'                 # T  Lan
$ws.Cells.Item($row, 21).Value = 'Varberg'                 # U  Kommun
$ws.Cells.Item($row, 22).Value = 'Halland'                 # V  Provins
$ws.Cells.Item($row, 23).Value = 'Rolfstorp'                # W  Forsamling

# Y3 / AA3 hold the literal text "2023-02-12" (not a real date value), so
# it also needs the quote-prefix trick to avoid becoming a date serial.
$ws.Cells.Item($row, 25).Value = "'" + '2023-02-12'        # Y  Startdatum
$ws.Cells.Item($row, 26).Value = '12:30'                   # Z  Starttid
$ws.Cells.Item($row, 27).Value = "'" + '2023-02-12'        # AA Slutdatum
$ws.Cells.Item($row, 28).Value = '12:30'                   # AB Sluttid

$ws.Cells.Item($row, 29).Value = 'Det pulserande typiska hoandet hördes precis norr om positionen.' # AC Publik kommentar

$ws.Cells.Item($row, 30).Value = $false                    # AD Ej aterfunnen
$ws.Cells.Item($row, 31).Value = $false                    # AE Osaker artbestamning
$ws.Cells.Item($row, 33).Value = $false                    # AG Ospontan

$ws.Cells.Item($row, 46).NumberFormat = "@"
$ws.Cells.Item($row, 46).Value = ""                        # AT Bestamningsar

$ws.Cells.Item($row, 49).Value = 'Andreas Källman'          # AW Rapportor
$ws.Cells.Item($row, 50).Value = 'Andreas Källman'          # AX Observatorer

$ws.Cells.Item($row, 51).NumberFormat = "@"
$ws.Cells.Item($row, 51).Value = ""                        # AY Projektnamn
